# Apply updated Price (D) and Volume(1h) (E) values to cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.004.07"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "2.756.44"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +0.03%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "578.24"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +0.06%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "157.95"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +2.96%  "
$ws.Range("E7").Value = "  -0.03%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.607"
$cell.ClearFormats()
$ws.Range("E8").Value = "  -0.14%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.111"
$cell.ClearFormats()
$ws.Range("E9").Value = "  -1.63%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "5.75"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -14.50%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.385"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("E12").Value = "  -2.65%  "
$ws.Range("D13").Value = "3.240.45"
$ws.Range("E13").Value = "  +0.40%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "26.89"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").Value = "63.680.05"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").Value = "2.757.45"
$ws.Range("E17").Value = "  +0.35%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "12.08"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +0.89%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "4.87"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +0.25%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "359.97"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -0.21%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.83"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -1.65%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.551"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +2.55%  "
$ws.Range("E23").Value = "  +0.18%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "65.80"
$cell.ClearFormats()
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("E26").Value = "  -0.04%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "8.48"
$cell.ClearFormats()
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").Value = "0.0₃0927"
$ws.Range("E28").Value = "  +2.79%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.95"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -2.53%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "7.00"
$cell.ClearFormats()
$ws.Range("E30").Value = "  -1.67%  "
$ws.Range("E31").Value = "  +1.08%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "168.20"
$cell.ClearFormats()
$ws.Range("E32").Value = "  -2.11%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "20.32"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -0.87%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "4.93"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +3.46%  "
$ws.Range("E35").Value = "  -0.04%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.46"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +1.22%  "
$ws.Range("E37").Value = "  -0.94%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.991"
$cell.ClearFormats()
$ws.Range("E38").Value = "  -0.21%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "6.16"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +11.95%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "4.15"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("E41").Value = "  -4.10%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "39.33"
$cell.ClearFormats()
$ws.Range("E42").Value = "  +0.62%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "21.54"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -1.01%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.0592"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +0.61%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "21.72"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -0.39%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.635"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("E47").Value = "  -0.20%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "135.66"
$cell.ClearFormats()
$ws.Range("E48").Value = "  -2.75%  "
$ws.Range("E49").Value = "  +0.16%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -0.03%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "11.04"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +0.56%  "
